# Update the "想去人数" (interest count) column F figures to the values
# captured at the newer scrape (commit "Update gh-pages to output
# generated at 456a3b4"). Only column F values change; everything else
# on each row (place, date, price, links, ...) stays as-is.

$wb = $excel.ActiveWorkbook

# Sheet 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2716
$ws1.Range("F7").Value = 2314
$ws1.Range("F8").Value = 1842
$ws1.Range("F9").Value = 218
$ws1.Range("F11").Value = 2487
$ws1.Range("F14").Value = 56
$ws1.Range("F16").Value = 130
$ws1.Range("F17").Value = 117
$ws1.Range("F18").Value = 9260
$ws1.Range("F20").Value = 7197
$ws1.Range("F21").Value = 11750
$ws1.Range("F24").Value = 234
$ws1.Range("F25").Value = 360
$ws1.Range("F26").Value = 563
$ws1.Range("F27").Value = 2606
$ws1.Range("F29").Value = 198
$ws1.Range("F30").Value = 2555
$ws1.Range("F31").Value = 728
$ws1.Range("F32").Value = 49
$ws1.Range("F33").Value = 4520
$ws1.Range("F34").Value = 934
$ws1.Range("F37").Value = 534

# Sheet 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 632
$ws3.Range("F4").Value = 160

# Sheet 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 632
$ws4.Range("F5").Value = 2716
$ws4.Range("F10").Value = 2314
$ws4.Range("F12").Value = 1842
$ws4.Range("F14").Value = 218
$ws4.Range("F15").Value = 2487
$ws4.Range("F19").Value = 56
$ws4.Range("F21").Value = 130
$ws4.Range("F22").Value = 117
$ws4.Range("F23").Value = 9261
$ws4.Range("F25").Value = 7197
$ws4.Range("F26").Value = 11750
$ws4.Range("F29").Value = 234
$ws4.Range("F30").Value = 360
$ws4.Range("F32").Value = 563
$ws4.Range("F34").Value = 2606
$ws4.Range("F38").Value = 198
$ws4.Range("F39").Value = 49
$ws4.Range("F40").Value = 4520
$ws4.Range("F45").Value = 535
